# Update cryptos list with latest prices and percentage changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.538.36'
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = '2.027.94'
$ws.Range("E3").Value = '  +1.57%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '255.94'
$ws.Range("E5").Value = '  +3.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.615'
$ws.Range("E6").Value = '  -2.67%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.11'
$ws.Range("E8").Value = '  -8.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.384'
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0786'
$ws.Range("E10").Value = '  -2.38%  '
$ws.Range("E11").Value = '  -2.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.53'
$ws.Range("E12").Value = '  -2.53%  '
$ws.Range("D13").Value = '2.321.29'
$ws.Range("E13").Value = '  +2.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.818'
$ws.Range("E14").Value = '  -3.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.21'
$ws.Range("E15").Value = '  -6.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.36'
$ws.Range("E16").Value = '  -1.68%  '
$ws.Range("D17").Value = '2.035.65'
$ws.Range("E17").Value = '  +2.59%  '
$ws.Range("D18").Value = '37.477.63'
$ws.Range("E18").Value = '  +0.85%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.59'
$ws.Range("E19").Value = '  -0.99%  '
$ws.Range("D20").Value = '0.0₃0850'
$ws.Range("E20").Value = '  -1.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.21'
$ws.Range("E21").Value = '  +0.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.23'
$ws.Range("E22").Value = '  -1.54%  '
$ws.Range("B23").Value = 'PancakeSwap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.63'
$ws.Range("E23").Value = '  +4.03%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("E25").Value = '  -1.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.04'
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.06'
$ws.Range("E27").Value = '  -2.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.91'
$ws.Range("E28").Value = '  +0.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.132'
$ws.Range("E29").Value = '  -9.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.38'
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.121'
$ws.Range("E31").Value = '  -1.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0668'
$ws.Range("E32").Value = '  +6.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.72'
$ws.Range("E33").Value = '  -3.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.58'
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.44'
$ws.Range("E35").Value = '  +5.19%  '
$ws.Range("E36").Value = '  +0.18%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.41'
$ws.Range("E37").Value = '  +1.37%  '
$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.82'
$ws.Range("E38").Value = '  +1.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.36'
$ws.Range("E39").Value = '  -2.75%  '
$ws.Range("E40").Value = '  +3.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0968'
$ws.Range("E41").Value = '  -1.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.20'
$ws.Range("E42").Value = '  +1.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0216'
$ws.Range("E43").Value = '  +0.81%  '
$ws.Range("D44").Value = '1.407.54'
$ws.Range("E44").Value = '  +1.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.03'
$ws.Range("E45").Value = '  -4.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.98'
$ws.Range("E46").Value = '  +0.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.04'
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.33'
$ws.Range("E48").Value = '  +1.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.87'
$ws.Range("E49").Value = '  +1.53%  '
$ws.Range("E50").Value = '  +1.66%  '
$ws.Range("D51").Value = '2.212.15'
$ws.Range("E51").Value = '  +2.13%  '
